$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sonderaufgaben")

$ws.Range("C3").Value = "Yavgaan"
$ws.Range("E3").Value = "Timera"

$ws.Range("C4").Value = "Kalbfleisch"

$ws.Range("C6").Value = "Üzülmez"

$ws.Range("E7").Value = "Delgado"

$ws.Range("C8").Value = "a.D"
$ws.Range("E8").Value = "a.D"

$ws.Range("E9").Value = "Moeeni Mahvelati"

$ws.Range("C10").Value = "Idic"
$ws.Range("E10").Value = "Rivola"

$wb.Save()
